# Updated cryptos list on Thu Jul  4 17:33:16 UTC 2024 with GitHub Actions
# Refresh Price (col D) and Volume(1h) (col E) figures for the crypto
# list, plus a ranking swap between EnergySwap and OKB (rows 42/43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.257.34'
$ws.Range('E2').Value = '  -3.54%  '

$ws.Range('D3').Value = '3.140.12'
$ws.Range('E3').Value = '  -5.08%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').Value = '524.47'
$ws.Range('E5').Value = '  -6.05%  '

$ws.Range('D6').Value = '135.27'
$ws.Range('E6').Value = '  -5.17%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '3.141.86'
$ws.Range('E8').Value = '  -5.05%  '

$ws.Range('E9').Value = '  -5.12%  '

$ws.Range('E10').Value = '  -7.24%  '

$ws.Range('E11').Value = '  -8.56%  '

$ws.Range('E12').Value = '  -6.32%  '

$ws.Range('D13').Value = '3.675.75'
$ws.Range('E13').Value = '  -5.17%  '

$ws.Range('E14').Value = '  -1.34%  '

$ws.Range('D15').Value = '25.61'
$ws.Range('E15').Value = '  -4.60%  '

$ws.Range('D16').Value = '3.137.26'
$ws.Range('E16').Value = '  -5.26%  '

$ws.Range('D17').Value = '58.204.44'
$ws.Range('E17').Value = '  -3.67%  '

$ws.Range('E18').Value = '  -7.43%  '

$ws.Range('D19').Value = '5.80'
$ws.Range('E19').Value = '  -5.26%  '

$ws.Range('D20').Value = '13.10'
$ws.Range('E20').Value = '  -7.69%  '

$ws.Range('D21').Value = '7.97'
$ws.Range('E21').Value = '  -7.98%  '

$ws.Range('D22').Value = '345.04'
$ws.Range('E22').Value = '  -8.06%  '

$ws.Range('E23').Value = '  +0.19%  '

$ws.Range('D24').Value = '68.71'
$ws.Range('E24').Value = '  -7.94%  '

$ws.Range('E25').Value = '  -5.46%  '

$ws.Range('D26').Value = '3.265.12'
$ws.Range('E26').Value = '  -5.27%  '

$ws.Range('E27').Value = '  -1.67%  '

$ws.Range('D28').Value = '0.0₃0960'
$ws.Range('E28').Value = '  -6.45%  '

$ws.Range('E29').Value = '  +0.30%  '

$ws.Range('E30').Value = '  -5.22%  '

$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.04%  '

$ws.Range('E32').Value = '  -8.83%  '

$ws.Range('E33').Value = '  -9.27%  '

$ws.Range('E34').Value = '  -4.80%  '

$ws.Range('D35').Value = '1.22'
$ws.Range('E35').Value = '  -1.66%  '

$ws.Range('D36').Value = '4.83'
$ws.Range('E36').Value = '  -5.85%  '

$ws.Range('D37').Value = '157.32'
$ws.Range('E37').Value = '  -5.56%  '

$ws.Range('D38').Value = '6.23'

$ws.Range('E39').Value = '  -10.37%  '

$ws.Range('D40').Value = '0.0693'
$ws.Range('E40').Value = '  -4.91%  '

$ws.Range('D41').Value = '3.171.10'
$ws.Range('E41').Value = '  -5.02%  '

# Rows 42/43 swap ranking order: OKB now above EnergySwap
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').Value = '40.49'
$ws.Range('E42').Value = '  -3.50%  '

$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '24.24'
$ws.Range('E43').Value = '  -9.46%  '

$ws.Range('D44').Value = '0.693'
$ws.Range('E44').Value = '  -7.77%  '

$ws.Range('E45').Value = '  -2.32%  '

$ws.Range('D46').Value = '3.91'
$ws.Range('E46').Value = '  -5.73%  '

$ws.Range('E47').Value = '  -0.01%  '

$ws.Range('E48').Value = '  -8.63%  '

$ws.Range('D49').Value = '2.262.92'
$ws.Range('E49').Value = '  -4.20%  '

$ws.Range('E50').Value = '  -3.37%  '

$ws.Range('D51').Value = '20.60'
$ws.Range('E51').Value = '  -3.27%  '
